$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 1152.5
$ws.Range("J10").Value = 1152.5
$ws.Range("L10").Value = 1152.5
$ws.Range("N10").Value = -1738.5
$ws.Range("H11").Value = 143
$ws.Range("I11").Value = 143
$ws.Range("K11").Value = 143
$ws.Range("M11").Value = -3
$ws.Range("H17").Value = 1552299.9
$ws.Range("I17").Value = 80
$ws.Range("J17").Value = 1607736.2
$ws.Range("K17").Value = 240
$ws.Range("L17").Value = 4823208.6
$ws.Range("M17").Value = -72
$ws.Range("N17").Value = -4823544.6
$ws.Range("H111").Value = 3250.1667
$ws.Range("I111").Value = 3225.65
$ws.Range("J111").Value = 3372.75
$ws.Range("K111").Value = 9676.950000000001
$ws.Range("L111").Value = 10118.25
$ws.Range("M111").Value = -6609.950000000001
$ws.Range("N111").Value = -16252.25
$ws.Range("H137").Value = 3774.795
$ws.Range("I137").Value = 1411.3334
$ws.Range("J137").Value = 7556.3335
$ws.Range("K137").Value = 4234.0002
$ws.Range("L137").Value = 22669.0005
$ws.Range("M137").Value = -1684.0002
$ws.Range("N137").Value = -27769.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 260
$ws.Range("I4").Value = 260
$ws.Range("K4").Value = 260
$ws.Range("M4").Value = -144
$ws.Range("H5").Value = 822.6923
$ws.Range("I5").Value = 975.75
$ws.Range("J5").Value = 577.8
$ws.Range("K5").Value = 975.75
$ws.Range("L5").Value = 577.8
$ws.Range("M5").Value = -863.75
$ws.Range("N5").Value = -801.8
$ws.Range("H28").Value = 6087.273
$ws.Range("I28").Value = 4710.1113
$ws.Range("K28").Value = 4710.1113
$ws.Range("M28").Value = -4518.1113
$ws.Range("H32").Value = 2470.759
$ws.Range("I32").Value = 1274.0541
$ws.Range("K32").Value = 1274.0541
$ws.Range("M32").Value = -987.0541000000001
$ws.Range("H41").Value = 4416.6924
$ws.Range("I41").Value = 4576.4165
$ws.Range("J41").Value = 2500
$ws.Range("K41").Value = 4576.4165
$ws.Range("L41").Value = 2500
$ws.Range("M41").Value = -4162.4165
$ws.Range("N41").Value = -3328
$ws.Range("H61").Value = 13890562
$ws.Range("I61").Value = 17242832
$ws.Range("K61").Value = 17242832
$ws.Range("M61").Value = -17242620
$ws.Range("H99").Value = 6087.273
$ws.Range("I99").Value = 4710.1113
$ws.Range("K99").Value = 4710.1113
$ws.Range("M99").Value = -1715.1113
$ws.Range("H132").Value = 28621446
$ws.Range("I132").Value = 11871.857
$ws.Range("K132").Value = 35615.571
$ws.Range("M132").Value = -33085.571
$ws.Range("H136").Value = 13890562
$ws.Range("I136").Value = 17242832
$ws.Range("K136").Value = 51728496
$ws.Range("M136").Value = -51725946

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 822.6923
$ws.Range("I4").Value = 975.75
$ws.Range("J4").Value = 577.8
$ws.Range("K4").Value = 975.75
$ws.Range("L4").Value = 577.8
$ws.Range("M4").Value = -860.75
$ws.Range("N4").Value = -807.8
$ws.Range("H19").Value = 3000
$ws.Range("I19").Value = 3000
$ws.Range("K19").Value = 3000
$ws.Range("M19").Value = -2827
$ws.Range("H99").Value = 4966.5
$ws.Range("I99").Value = 4514.143
$ws.Range("K99").Value = 4514.143
$ws.Range("M99").Value = -3016.143
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("M112").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1476.7142
$ws.Range("I16").Value = 862.4
$ws.Range("K16").Value = 862.4
$ws.Range("M16").Value = -575.4
$ws.Range("H31").Value = 5105578
$ws.Range("I31").Value = 1675.6111
$ws.Range("J31").Value = 19239462
$ws.Range("K31").Value = 1675.6111
$ws.Range("L31").Value = 19239462
$ws.Range("M31").Value = -1380.6111
$ws.Range("N31").Value = -19240052
$ws.Range("H34").Value = 5105578
$ws.Range("I34").Value = 1675.6111
$ws.Range("J34").Value = 19239462
$ws.Range("K34").Value = 1675.6111
$ws.Range("L34").Value = 19239462
$ws.Range("M34").Value = -1473.6111
$ws.Range("N34").Value = -19239866
$ws.Range("H43").Value = 14628.333
$ws.Range("J43").Value = 14628.333
$ws.Range("L43").Value = 14628.333
$ws.Range("N43").Value = -14996.333
$ws.Range("H58").Value = 3029.1292
$ws.Range("I58").Value = 3539
$ws.Range("J58").Value = 2102.0908
$ws.Range("K58").Value = 3539
$ws.Range("L58").Value = 2102.0908
$ws.Range("M58").Value = -3336
$ws.Range("N58").Value = -2508.0908
$ws.Range("H93").Value = 17639.234
$ws.Range("I93").Value = 13191.2
$ws.Range("K93").Value = 13191.2
$ws.Range("M93").Value = -11319.2
$ws.Range("H101").Value = 14628.333
$ws.Range("J101").Value = 14628.333
$ws.Range("L101").Value = 14628.333
$ws.Range("N101").Value = -21118.333
$ws.Range("H113").Value = 1476.7142
$ws.Range("I113").Value = 862.4
$ws.Range("K113").Value = 862.4
$ws.Range("M113").Value = 1307.6
$ws.Range("H122").Value = 4387314
$ws.Range("I122").Value = 1260.1111
$ws.Range("J122").Value = 17545476
$ws.Range("K122").Value = 3780.3333
$ws.Range("L122").Value = 52636428
$ws.Range("M122").Value = -1330.3333
$ws.Range("N122").Value = -52641328
$ws.Range("H136").Value = 3029.1292
$ws.Range("I136").Value = 3539
$ws.Range("J136").Value = 2102.0908
$ws.Range("K136").Value = 10617
$ws.Range("L136").Value = 6306.2724
$ws.Range("M136").Value = -8067
$ws.Range("N136").Value = -11406.2724

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1166.6888
$ws.Range("J2").Value = 54.358974
$ws.Range("L2").Value = 326.153844
$ws.Range("N2").Value = -552.153844
$ws.Range("H12").Value = 446.89474
$ws.Range("J12").Value = 527.7143
$ws.Range("L12").Value = 1583.1429
$ws.Range("N12").Value = -1929.1429
$ws.Range("H46").Value = 348
$ws.Range("I46").Value = 22.5
$ws.Range("K46").Value = 67.5
$ws.Range("M46").Value = 23.5
$ws.Range("H50").Value = 1195
$ws.Range("I50").Value = 725.8333
$ws.Range("J50").Value = 2133.3333
$ws.Range("K50").Value = 2177.4999
$ws.Range("L50").Value = 6399.999899999999
$ws.Range("M50").Value = -1696.4999
$ws.Range("N50").Value = -7361.999899999999
$ws.Range("H53").Value = 1195
$ws.Range("I53").Value = 725.8333
$ws.Range("J53").Value = 2133.3333
$ws.Range("K53").Value = 2177.4999
$ws.Range("L53").Value = 6399.999899999999
$ws.Range("M53").Value = -1696.4999
$ws.Range("N53").Value = -7361.999899999999
$ws.Range("H55").Value = 5824066.5
$ws.Range("J55").Value = 9529395
$ws.Range("L55").Value = 28588185
$ws.Range("N55").Value = -28588539
$ws.Range("H68").Value = 2813.2307
$ws.Range("J68").Value = 2957.6296
$ws.Range("L68").Value = 8872.888800000001
$ws.Range("N68").Value = -10494.8888
$ws.Range("H71").Value = 2813.2307
$ws.Range("J71").Value = 2957.6296
$ws.Range("L71").Value = 26618.6664
$ws.Range("N71").Value = -34730.6664
$ws.Range("H107").Value = 1536.4546
$ws.Range("J107").Value = 1894
$ws.Range("L107").Value = 5682
$ws.Range("N107").Value = -9522

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3408.75
$ws.Range("I80").Value = 3181.4285
$ws.Range("K80").Value = 3181.4285
$ws.Range("M80").Value = -2183.4285
$ws.Range("H83").Value = 3408.75
$ws.Range("I83").Value = 3181.4285
$ws.Range("K83").Value = 15907.1425
$ws.Range("M83").Value = -10915.1425

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2926.9614
$ws.Range("I22").Value = 2512.889
$ws.Range("J22").Value = 3858.625
$ws.Range("K22").Value = 2512.889
$ws.Range("L22").Value = 3858.625
$ws.Range("M22").Value = -2217.889
$ws.Range("N22").Value = -4448.625
$ws.Range("H27").Value = 2926.9614
$ws.Range("I27").Value = 2512.889
$ws.Range("J27").Value = 3858.625
$ws.Range("K27").Value = 2512.889
$ws.Range("L27").Value = 3858.625
$ws.Range("M27").Value = -2405.889
$ws.Range("N27").Value = -4072.625
$ws.Range("H46").Value = 1795.5938
$ws.Range("J46").Value = 5416.6665
$ws.Range("L46").Value = 5416.6665
$ws.Range("N46").Value = -5792.6665
$ws.Range("H54").Value = 33000
$ws.Range("J54").Value = 33000
$ws.Range("L54").Value = 33000
$ws.Range("N54").Value = -34288
$ws.Range("H55").Value = 560.4483
$ws.Range("I55").Value = 233.61905
$ws.Range("J55").Value = 1418.375
$ws.Range("K55").Value = 233.61905
$ws.Range("L55").Value = 1418.375
$ws.Range("M55").Value = -60.61904999999999
$ws.Range("N55").Value = -1764.375
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H140").Value = 75388.336
$ws.Range("J140").Value = 75385
$ws.Range("L140").Value = 75385
$ws.Range("N140").Value = -85745

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2953.3
$ws.Range("I107").Value = 2522.1667
$ws.Range("J107").Value = 3600
$ws.Range("K107").Value = 7566.500100000001
$ws.Range("L107").Value = 10800
$ws.Range("M107").Value = -5646.500100000001
$ws.Range("N107").Value = -14640
$ws.Range("H135").Value = 50047856
$ws.Range("J135").Value = 50047856
$ws.Range("L135").Value = 50047856
$ws.Range("N135").Value = -50057996
